# Slide 2: "Chirp - ChiptuneSAK Intermediate Representation" diagram.
# Reposition/resize the big background rectangle and two labels to make
# room for the newly added Chirp-transformation documentation.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# "Rectangle 2": background rectangle behind the Chirp diagram -- moved
# down a bit and made shorter (off x/y 685800/1752600 -> 685800/2121932;
# ext cx/cy 7772400/3352800 -> 7772400/2983468).
$rect = $s.Shapes.Item("Rectangle 2")
$rect.Top = 167.08126068115234
$rect.Height = 234.91873931884768

# "TextBox 3": the "ChirpSong" label -- nudged down
# (off x/y 3886200/1905000 -> 3886200/2133600).
$chirpSongLabel = $s.Shapes.Item("TextBox 3")
$chirpSongLabel.Top = 167.99999237060547

# "TextBox 12": the "metadata" label -- moved right and down
# (off x/y 3953119/2362200 -> 3977349/2435423).
$metadataLabel = $s.Shapes.Item("TextBox 12")
$metadataLabel.Left = 313.17707824707037
$metadataLabel.Top = 191.76558685302734
